$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update alternative flow numbering: "Alternativa 1" now refers to step 5 (was step 6)
$ws.Range("B14").Value = "Alternativa 1`n[Cliente já registado]`n(Passo 5)"
$ws.Range("D14").Value = "5.1. Informa cliente que já está registado"
$ws.Range("D15").Value = "5.2. Sai do ecrã de registo"

# D17 belongs to "Alternativa 2" block and its numbering shifts from 6.1.1 to 5.1.1
$ws.Range("D17").Value = "5.1.1. Informa cliente que os dados são inválidos"

# Re-assigning the multi-line text to the merged B14 cell makes Excel recompute
# the row's autofit height; restore the original row height afterwards.
$ws.Rows.Item(14).RowHeight = 19.5

$ws.Range("F12").Select()

$wb.Save()
